$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032896927620757
$ws.Range("D2").Value = 1.033823092741411
$ws.Range("E2").Value = 1.032220749832632
$ws.Range("F2").Value = 1.031523158370919
$ws.Range("I2").Value = 1.031697688990885
$ws.Range("J2").Value = 1.038024625739814
$ws.Range("K2").Value = 1.036624334470686
$ws.Range("L2").Value = 1.035026607613834
$ws.Range("M2").Value = 1.034331030652835
$ws.Range("N2").Value = 1.039498739771474
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034486855960751
$ws.Range("D3").Value = 1.035240211853633
$ws.Range("E3").Value = 1.033591638336064
$ws.Range("F3").Value = 1.033752985904857
$ws.Range("I3").Value = 1.03198892692478
$ws.Range("J3").Value = 1.039253792883857
$ws.Range("K3").Value = 1.037848892624037
$ws.Range("L3").Value = 1.036204719180097
$ws.Range("M3").Value = 1.036365635439297
$ws.Range("N3").Value = 1.040729652473829
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035512635041849
$ws.Range("D4").Value = 1.036154628480346
$ws.Range("E4").Value = 1.034476299611775
$ws.Range("F4").Value = 1.035192246323551
$ws.Range("I4").Value = 1.032174825122325
$ws.Range("J4").Value = 1.04004582817573
$ws.Range("K4").Value = 1.038638202201568
$ws.Range("L4").Value = 1.036964123708643
$ws.Range("M4").Value = 1.03767825546546
$ws.Range("N4").Value = 1.041522812546663
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035943165628421
$ws.Range("D5").Value = 1.036538449272424
$ws.Range("E5").Value = 1.034847648685751
$ws.Range("F5").Value = 1.035796477591118
$ws.Range("I5").Value = 1.032252368926881
$ws.Range("J5").Value = 1.040378016550308
$ws.Range("K5").Value = 1.038969305765757
$ws.Range("L5").Value = 1.037282690312303
$ws.Range("M5").Value = 1.038229168830073
$ws.Range("N5").Value = 1.041855472666836
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036015412553412
$ws.Range("D6").Value = 1.036602859622997
$ws.Range("E6").Value = 1.034909967157071
$ws.Range("F6").Value = 1.035897882504534
$ws.Range("I6").Value = 1.032265353323789
$ws.Range("J6").Value = 1.040433746829015
$ws.Range("K6").Value = 1.039024857471601
$ws.Range("L6").Value = 1.037336139051051
$ws.Range("M6").Value = 1.038321616781854
$ws.Range("N6").Value = 1.041911282088931
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035518390571422
$ws.Range("D7").Value = 1.036159759454168
$ws.Range("E7").Value = 1.034481263791649
$ws.Range("F7").Value = 1.03520032333385
$ws.Range("I7").Value = 1.032175863650222
$ws.Range("J7").Value = 1.040050269955541
$ws.Range("K7").Value = 1.038642629246564
$ws.Range("L7").Value = 1.036968383096899
$ws.Range("M7").Value = 1.037685620349742
$ws.Range("N7").Value = 1.041527260634311
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033434883156682
$ws.Range("D8").Value = 1.034302550614932
$ws.Range("E8").Value = 1.032684550570059
$ws.Range("F8").Value = 1.032277496847753
$ws.Range("I8").Value = 1.031796644611757
$ws.Range("J8").Value = 1.038440722259841
$ws.Range("K8").Value = 1.037038820013962
$ws.Range("L8").Value = 1.035425365132631
$ws.Range("M8").Value = 1.035019457308539
$ws.Range("N8").Value = 1.039915427196266
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029739791062834
$ws.Range("D9").Value = 1.031009828712827
$ws.Range("E9").Value = 1.029499674961994
$ws.Range("F9").Value = 1.027098504391734
$ws.Range("I9").Value = 1.03110871982675
$ws.Range("J9").Value = 1.035578580505082
$ws.Range("K9").Value = 1.034188771465766
$ws.Range("L9").Value = 1.032683604473757
$ws.Range("M9").Value = 1.03029039606787
$ws.Range("N9").Value = 1.037049220871923
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027259590143533
$ws.Range("D10").Value = 1.028800445113052
$ws.Range("E10").Value = 1.027363065406554
$ws.Range("F10").Value = 1.023625032571047
$ws.Range("I10").Value = 1.030636657430292
$ws.Range("J10").Value = 1.033652357928665
$ws.Range("K10").Value = 1.032271965615746
$ws.Range("L10").Value = 1.03083980923789
$ws.Range("M10").Value = 1.027115436037157
$ws.Range("N10").Value = 1.035120262838506
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02618145185855
$ws.Range("D11").Value = 1.027840216392136
$ws.Range("E11").Value = 1.026434565325906
$ws.Range("F11").Value = 1.022115690811918
$ws.Range("I11").Value = 1.030429013475902
$ws.Range("J11").Value = 1.032813824004838
$ws.Range("K11").Value = 1.031437835572051
$ws.Range("L11").Value = 1.030037495573137
$ws.Range("M11").Value = 1.025735040680873
$ws.Range("N11").Value = 1.034280538100326
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025780335890226
$ws.Range("D12").Value = 1.027482996865388
$ws.Range("E12").Value = 1.026089164589993
$ws.Range("F12").Value = 1.021554226332234
$ws.Range("I12").Value = 1.030351394726526
$ws.Range("J12").Value = 1.032501670282798
$ws.Range("K12").Value = 1.031127367120728
$ws.Range("L12").Value = 1.029738876210821
$ws.Range("M12").Value = 1.025221428432523
$ws.Range("N12").Value = 1.033967941084196
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025866406201445
$ws.Range("D13").Value = 1.027559646688693
$ws.Range("E13").Value = 1.02616327772059
$ws.Range("F13").Value = 1.021674700322787
$ws.Range("I13").Value = 1.030368066476042
$ws.Range("J13").Value = 1.032568659531819
$ws.Range("K13").Value = 1.03119399262756
$ws.Range("L13").Value = 1.029802958630948
$ws.Range("M13").Value = 1.025331639916092
$ws.Range("N13").Value = 1.034035025465635
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026148308794454
$ws.Range("D14").Value = 1.027810699773633
$ws.Range("E14").Value = 1.026406024967355
$ws.Range("F14").Value = 1.022069297047469
$ws.Range("I14").Value = 1.030422607516617
$ws.Range("J14").Value = 1.032788035332759
$ws.Range("K14").Value = 1.031412185193166
$ws.Range("L14").Value = 1.0300128240059
$ws.Range("M14").Value = 1.025692603283687
$ws.Range("N14").Value = 1.034254712805375
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026321911981133
$ws.Range("D15").Value = 1.02796530884112
$ws.Range("E15").Value = 1.026555520977375
$ws.Range("F15").Value = 1.022312310504311
$ws.Range("I15").Value = 1.030456146927075
$ws.Range("J15").Value = 1.03292310883396
$ws.Range("K15").Value = 1.031546536233004
$ws.Range("L15").Value = 1.03014204852609
$ws.Range("M15").Value = 1.025914888338498
$ws.Range("N15").Value = 1.034389978126441
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027331052788458
$ws.Range("D16").Value = 1.028864096296344
$ws.Range("E16").Value = 1.027424615501016
$ws.Range("F16").Value = 1.023725088224395
$ws.Range("I16").Value = 1.030650369510496
$ws.Range("J16").Value = 1.033707913379139
$ws.Range("K16").Value = 1.032327235708914
$ws.Range("L16").Value = 1.030892972102508
$ws.Range("M16").Value = 1.027206927602908
$ws.Range("N16").Value = 1.035175897184093
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027962925007426
$ws.Range("D17").Value = 1.02942692115532
$ws.Range("E17").Value = 1.027968873398426
$ws.Range("F17").Value = 1.024609844411698
$ws.Range("I17").Value = 1.030771330478101
$ws.Range("J17").Value = 1.034198994955858
$ws.Range("K17").Value = 1.032815829972309
$ws.Range("L17").Value = 1.031362943366914
$ws.Range("M17").Value = 1.028015866223682
$ws.Range("N17").Value = 1.035667676152978
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02833108180421
$ws.Range("D18").Value = 1.029754865726501
$ws.Range("E18").Value = 1.028286009060539
$ws.Range("F18").Value = 1.025125397909577
$ws.Range("I18").Value = 1.030841572840976
$ws.Range("J18").Value = 1.034485004487535
$ws.Range("K18").Value = 1.033100420156731
$ws.Range("L18").Value = 1.031636690544158
$ws.Range("M18").Value = 1.028487166420562
$ws.Range("N18").Value = 1.035954091850996
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028456545820011
$ws.Range("D19").Value = 1.02986662883154
$ws.Range("E19").Value = 1.02839409012661
$ws.Range("F19").Value = 1.025301102776804
$ws.Range("I19").Value = 1.030865470852857
$ws.Range("J19").Value = 1.034582453818345
$ws.Range("K19").Value = 1.033197390849853
$ws.Range("L19").Value = 1.031729967334691
$ws.Range("M19").Value = 1.028647776857222
$ws.Range("N19").Value = 1.036051679571033
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027895172966087
$ws.Range("D20").Value = 1.029366570815903
$ws.Range("E20").Value = 1.027910512952512
$ws.Range("F20").Value = 1.024514971399633
$ws.Range("I20").Value = 1.030758384824911
$ws.Range("J20").Value = 1.034146351131131
$ws.Range("K20").Value = 1.032763449755491
$ws.Range("L20").Value = 1.031312559202582
$ws.Range("M20").Value = 1.027929130824086
$ws.Range("N20").Value = 1.035614957567982
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026065313517561
$ws.Range("D21").Value = 1.02773678613608
$ws.Range("E21").Value = 1.026334556265051
$ws.Range("F21").Value = 1.021953121245174
$ws.Range("I21").Value = 1.030406560114026
$ws.Range("J21").Value = 1.03272345364512
$ws.Range("K21").Value = 1.031347950569571
$ws.Range("L21").Value = 1.029951040669268
$ws.Range("M21").Value = 1.025586332866391
$ws.Range("N21").Value = 1.034190039404332
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024911053903934
$ws.Range("D22").Value = 1.026708899039037
$ws.Range("E22").Value = 1.025340706363554
$ws.Range("F22").Value = 1.0203375798012
$ws.Range("I22").Value = 1.030182513773406
$ws.Range("J22").Value = 1.031824852007701
$ws.Range("K22").Value = 1.030454286997939
$ws.Range("L22").Value = 1.029091496974644
$ws.Range("M22").Value = 1.024108263768596
$ws.Range("N22").Value = 1.033290161649517
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025523310422205
$ws.Range("D23").Value = 1.027254107390533
$ws.Range("E23").Value = 1.025867852365514
$ws.Range("F23").Value = 1.021194474646518
$ws.Range("I23").Value = 1.030301555576175
$ws.Range("J23").Value = 1.032301598655534
$ws.Range("K23").Value = 1.030928388582179
$ws.Range("L23").Value = 1.029547493667373
$ws.Range("M23").Value = 1.024892305183155
$ws.Range("N23").Value = 1.033767585332272
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027925788443985
$ws.Range("D24").Value = 1.029393841592038
$ws.Range("E24").Value = 1.027936884516917
$ws.Range("F24").Value = 1.024557842004432
$ws.Range("I24").Value = 1.030764235372756
$ws.Range("J24").Value = 1.034170139938107
$ws.Range("K24").Value = 1.032787119352208
$ws.Range("L24").Value = 1.031335326808947
$ws.Range("M24").Value = 1.027968324485492
$ws.Range("N24").Value = 1.035638780157792
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030697957885978
$ws.Range("D25").Value = 1.031863531061262
$ws.Range("E25").Value = 1.030325344872316
$ws.Range("F25").Value = 1.028440948832606
$ws.Range("I25").Value = 1.031288919525301
$ws.Range("J25").Value = 1.036321655667884
$ws.Range("K25").Value = 1.034928483155612
$ws.Range("L25").Value = 1.033395179426279
$ws.Range("M25").Value = 1.031516786553059
$ws.Range("N25").Value = 1.037793351286687
